$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.409.42"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "1.629.26"
$ws.Range("E3").Value = "  -0.63%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3754"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3637"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08199"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.230"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.530"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001246"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.337"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.47%  "

$ws.Range("D17").Value = "1.629.47"
$ws.Range("E17").Value = "  -0.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06966"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.526"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.29%  "

$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("E23").Value = "  -1.20%  "

$ws.Range("D24").Value = "23.403.76"
$ws.Range("E24").Value = "  -0.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.131"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.456"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.295"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.36%  "

$ws.Range("D31").Value = "1.810.36"
$ws.Range("E31").Value = "  -0.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.236"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.794"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.034"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02781"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2514"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08771"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07107"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7015"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.29%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.342"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6536"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.288"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.973"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07997"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.202"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "126.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.17%  "
